# Insert a new weekly price group (3 rows: Especial / Primera / Segunda)
# for "Terminal La Palmera de La Serena - Frutilla" just above the existing
# row 448, shifting the subsequent rows (old 448-462) down to 451-465.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at row 448 (pushes everything below down by 3).
$ws.Rows.Item(448).Insert()
$ws.Rows.Item(448).Insert()
$ws.Rows.Item(448).Insert()

# Common values shared by the three new rows.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$fecha = 44509
$codreg = 4
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "`$/bandeja 7 kilos"
$origen = "Provincia de Melipilla"
$kgUnidad = 7

# Row 448: Especial
$ws.Range("A448").Value = $mercadoId
$ws.Range("B448").Value = $mercado
$ws.Range("C448").Value = $region
$ws.Range("D448").Value = $fecha
$ws.Range("E448").Value = $codreg
$ws.Range("F448").Value = $tipo
$ws.Range("G448").Value = $productoId
$ws.Range("H448").Value = $producto
$ws.Range("I448").Value = $categoriaId
$ws.Range("J448").Value = $categoria
$ws.Range("K448").Value = $variedad
$ws.Range("L448").Value = "Especial"
$ws.Range("M448").Value = 400
$ws.Range("N448").Value = 12500
$ws.Range("O448").Value = 13000
$ws.Range("P448").Value = 12750
$ws.Range("Q448").Value = $unidad
$ws.Range("R448").Value = $origen
$ws.Range("S448").Value = 1821
$ws.Range("T448").Value = $kgUnidad

# Row 449: Primera
$ws.Range("A449").Value = $mercadoId
$ws.Range("B449").Value = $mercado
$ws.Range("C449").Value = $region
$ws.Range("D449").Value = $fecha
$ws.Range("E449").Value = $codreg
$ws.Range("F449").Value = $tipo
$ws.Range("G449").Value = $productoId
$ws.Range("H449").Value = $producto
$ws.Range("I449").Value = $categoriaId
$ws.Range("J449").Value = $categoria
$ws.Range("K449").Value = $variedad
$ws.Range("L449").Value = "Primera"
$ws.Range("M449").Value = 320
$ws.Range("N449").Value = 10500
$ws.Range("O449").Value = 11000
$ws.Range("P449").Value = 10750
$ws.Range("Q449").Value = $unidad
$ws.Range("R449").Value = $origen
$ws.Range("S449").Value = 1536
$ws.Range("T449").Value = $kgUnidad

# Row 450: Segunda
$ws.Range("A450").Value = $mercadoId
$ws.Range("B450").Value = $mercado
$ws.Range("C450").Value = $region
$ws.Range("D450").Value = $fecha
$ws.Range("E450").Value = $codreg
$ws.Range("F450").Value = $tipo
$ws.Range("G450").Value = $productoId
$ws.Range("H450").Value = $producto
$ws.Range("I450").Value = $categoriaId
$ws.Range("J450").Value = $categoria
$ws.Range("K450").Value = $variedad
$ws.Range("L450").Value = "Segunda"
$ws.Range("M450").Value = 260
$ws.Range("N450").Value = 8500
$ws.Range("O450").Value = 9000
$ws.Range("P450").Value = 8750
$ws.Range("Q450").Value = $unidad
$ws.Range("R450").Value = $origen
$ws.Range("S450").Value = 1250
$ws.Range("T450").Value = $kgUnidad
